# Apply "New range for cryo test" edit:
# - Rows 34-37, columns B:BO set to 2500
# - Update frozen pane top-left cell and active selection on the Biases sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Biases")

# Set the new values for rows 34-37, columns B through BO
$ws.Range("B34:BO37").Value = 2500

# Update the view: frozen-pane scroll position and current selection
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 49   # column AW
$ws.Application.ActiveWindow.ScrollRow = 23

$ws.Range("B34:AE37").Select()
$ws.Application.ActiveWindow.RangeSelection.Item(1).Activate()
$ws.Cells.Item(34, 29).Activate()
